$wb = $excel.ActiveWorkbook

# --- Sheet1: fixed recourse summary table (corrected error in fixed recourse data) ---
$ws1 = $wb.Worksheets.Item("Sheet1")

$ws1.Range("B2").Value = -274.9655817221359
$ws1.Range("C2").Value = 28.361499513
$ws1.Range("F2").Value = 50.0
$ws1.Range("G2").Value = 27700.0
$ws1.Range("H2").Value = 30250.0
$ws1.Range("I2").Value = 2500.0

$ws1.Range("B3").Value = -273.9816225888468
$ws1.Range("C3").Value = 19.684519199
$ws1.Range("F3").Value = 50.0
$ws1.Range("G3").Value = 27700.0
$ws1.Range("H3").Value = 30250.0
$ws1.Range("I3").Value = 2500.0

$ws1.Range("B4").Value = -274.08960459636427
$ws1.Range("C4").Value = 14.988418612
$ws1.Range("F4").Value = 50.0
$ws1.Range("G4").Value = 27700.0
$ws1.Range("H4").Value = 30250.0
$ws1.Range("I4").Value = 2500.0

$ws1.Range("B5").Value = -276.86855154162515
$ws1.Range("C5").Value = 16.164857087
$ws1.Range("F5").Value = 50.0
$ws1.Range("G5").Value = 27700.0
$ws1.Range("H5").Value = 30250.0
$ws1.Range("I5").Value = 2500.0

$ws1.Range("B6").Value = -272.1653938311721
$ws1.Range("C6").Value = 17.695829122
$ws1.Range("F6").Value = 50.0
$ws1.Range("G6").Value = 27700.0
$ws1.Range("H6").Value = 30250.0
$ws1.Range("I6").Value = 2500.0

$ws1.Range("B7").Value = -268.97221187709164
$ws1.Range("C7").Value = 17.885493201
$ws1.Range("F7").Value = 50.0
$ws1.Range("G7").Value = 27700.0
$ws1.Range("H7").Value = 30250.0
$ws1.Range("I7").Value = 2500.0

$ws1.Range("B8").Value = -265.4281513734784
$ws1.Range("C8").Value = 19.47060509
$ws1.Range("F8").Value = 50.0
$ws1.Range("G8").Value = 27700.0
$ws1.Range("H8").Value = 30250.0
$ws1.Range("I8").Value = 2500.0

$ws1.Range("B9").Value = -274.2017084750978
$ws1.Range("C9").Value = 17.878042283
$ws1.Range("F9").Value = 50.0
$ws1.Range("G9").Value = 27700.0
$ws1.Range("H9").Value = 30250.0
$ws1.Range("I9").Value = 2500.0

$ws1.Range("B10").Value = -271.53604073578464
$ws1.Range("C10").Value = 18.787386758
$ws1.Range("F10").Value = 50.0
$ws1.Range("G10").Value = 27700.0
$ws1.Range("H10").Value = 30250.0
$ws1.Range("I10").Value = 2500.0

$ws1.Range("B11").Value = -268.78676349663795
$ws1.Range("C11").Value = 18.493577767
$ws1.Range("F11").Value = 50.0
$ws1.Range("G11").Value = 27700.0
$ws1.Range("H11").Value = 30250.0
$ws1.Range("I11").Value = 2500.0

# --- Per-instance MP (master problem) iteration detail sheets: changed MP time limit ---
$wsD = $wb.Worksheets.Item("1")
$wsD.Range("D2").Value = 0.9088385056567383
$wsD.Range("E2").Value = 70.77445
$wsD.Range("B3").Value = -274.9655817221359
$wsD.Range("C3").Value = 0.08770811021076466
$wsD.Range("D3").Value = 16.902787455143798

$wsD = $wb.Worksheets.Item("2")
$wsD.Range("D2").Value = 0.025314165014038086
$wsD.Range("E2").Value = 68.31078
$wsD.Range("B3").Value = -273.9816225888468
$wsD.Range("C3").Value = 0.09158070790911396
$wsD.Range("D3").Value = 18.891956181888673

$wsD = $wb.Worksheets.Item("3")
$wsD.Range("D2").Value = 0.05396996686218262
$wsD.Range("E2").Value = 72.03609
$wsD.Range("B3").Value = -274.08960459636427
$wsD.Range("C3").Value = 0.09384263339934619
$wsD.Range("D3").Value = 14.110434908917236

$wsD = $wb.Worksheets.Item("4")
$wsD.Range("D2").Value = 0.057115765817260744
$wsD.Range("E2").Value = 71.29895
$wsD.Range("B3").Value = -276.86855154162515
$wsD.Range("C3").Value = 0.0
$wsD.Range("D3").Value = 15.265561777960693

$wsD = $wb.Worksheets.Item("5")
$wsD.Range("D2").Value = 0.05462881648291015
$wsD.Range("E2").Value = 70.50756
$wsD.Range("B3").Value = -272.1653938311721
$wsD.Range("C3").Value = 0.08075643941171622
$wsD.Range("D3").Value = 16.889026603522094

$wsD = $wb.Worksheets.Item("6")
$wsD.Range("D2").Value = 0.052927086459960936
$wsD.Range("E2").Value = 74.8484
$wsD.Range("B3").Value = -268.97221187709164
$wsD.Range("C3").Value = 0.08639506344671546
$wsD.Range("D3").Value = 17.075744087620727

$wsD = $wb.Worksheets.Item("7")
$wsD.Range("D2").Value = 0.05447040664355469
$wsD.Range("E2").Value = 69.60139
$wsD.Range("B3").Value = -265.4281513734784
$wsD.Range("C3").Value = 0.06623651690862611
$wsD.Range("D3").Value = 18.586615919048096

$wsD = $wb.Worksheets.Item("8")
$wsD.Range("D2").Value = 0.049114984947387696
$wsD.Range("E2").Value = 71.06337
$wsD.Range("B3").Value = -274.2017084750978
$wsD.Range("C3").Value = 0.06405499379048205
$wsD.Range("D3").Value = 17.09294442499768

$wsD = $wb.Worksheets.Item("9")
$wsD.Range("D2").Value = 0.040193955549560546
$wsD.Range("E2").Value = 68.71954
$wsD.Range("B3").Value = -271.53604073578464
$wsD.Range("C3").Value = 0.07011464592711199
$wsD.Range("D3").Value = 17.900044588735472

$wsD = $wb.Worksheets.Item("10")
$wsD.Range("D2").Value = 0.0532828747911377
$wsD.Range("E2").Value = 70.3188
$wsD.Range("B3").Value = -268.78676349663795
$wsD.Range("C3").Value = 0.07971176734353576
$wsD.Range("D3").Value = 17.5638273662312
